$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 900
$ws.Range("I51").Value = 800
$ws.Range("J51").Value = 1000
$ws.Range("K51").Value = 800
$ws.Range("L51").Value = 1000
$ws.Range("M51").Value = -316
$ws.Range("N51").Value = -1968
$ws.Range("H63").Value = 597191.3
$ws.Range("J63").Value = 597191.3
$ws.Range("L63").Value = 597191.3
$ws.Range("N63").Value = -598439.3
$ws.Range("H66").Value = 597191.3
$ws.Range("J66").Value = 597191.3
$ws.Range("L66").Value = 1791573.9
$ws.Range("N66").Value = -1797813.9
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H95").Value = 42812
$ws.Range("J95").Value = 42812
$ws.Range("L95").Value = 42812
$ws.Range("N95").Value = -48304
$ws.Range("H96").Value = 587.6667
$ws.Range("I96").Value = 508.66666
$ws.Range("J96").Value = 666.6667
$ws.Range("K96").Value = 1525.99998
$ws.Range("L96").Value = 2000.0001
$ws.Range("M96").Value = -152.9999800000001
$ws.Range("N96").Value = -4746.0001
$ws.Range("H107").Value = 1748657.9
$ws.Range("I107").Value = 1923453.4
$ws.Range("J107").Value = 703
$ws.Range("K107").Value = 1923453.4
$ws.Range("L107").Value = 703
$ws.Range("M107").Value = -1921533.4
$ws.Range("N107").Value = -4543
$ws.Range("H108").Value = 39888
$ws.Range("J108").Value = 39888
$ws.Range("L108").Value = 39888
$ws.Range("N108").Value = -47568
$ws.Range("H116").Value = 4504.2573
$ws.Range("I116").Value = 4531
$ws.Range("J116").Value = 4475.9414
$ws.Range("K116").Value = 4531
$ws.Range("L116").Value = 4475.9414
$ws.Range("M116").Value = -1089
$ws.Range("N116").Value = -11359.9414
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H124").Value = 10709
$ws.Range("I124").Value = 10709
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 10709
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -5799
$ws.Range("N124").ClearContents()
$ws.Range("H135").Value = 22289.49
$ws.Range("I135").Value = 25305.072
$ws.Range("J135").Value = 1683
$ws.Range("K135").Value = 227745.648
$ws.Range("L135").Value = 15147
$ws.Range("M135").Value = -225210.648
$ws.Range("N135").Value = -20217
$ws.Range("H137").Value = 3001116
$ws.Range("I137").Value = 1389989.1
$ws.Range("J137").Value = 7144014
$ws.Range("K137").Value = 4169967.3
$ws.Range("L137").Value = 21432042
$ws.Range("M137").Value = -4167417.3
$ws.Range("N137").Value = -21437142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1240.5745
$ws.Range("I61").Value = 1333.7646
$ws.Range("J61").Value = 996.8461
$ws.Range("K61").Value = 1333.7646
$ws.Range("L61").Value = 996.8461
$ws.Range("M61").Value = -1121.7646
$ws.Range("N61").Value = -1420.8461
$ws.Range("H74").Value = 996.0857
$ws.Range("I74").Value = 989.7879
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 989.7879
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = -115.7879
$ws.Range("N74").Value = -2848
$ws.Range("H77").Value = 996.0857
$ws.Range("I77").Value = 989.7879
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 4948.9395
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = -580.9395000000004
$ws.Range("N77").Value = -14236
$ws.Range("H122").Value = 1666.4193
$ws.Range("I122").Value = 1559.25
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 4677.75
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -2227.75
$ws.Range("N122").Value = -12900.0001
$ws.Range("H132").Value = 110778.914
$ws.Range("I132").Value = 132984.38
$ws.Range("J132").Value = 5303
$ws.Range("K132").Value = 398953.14
$ws.Range("L132").Value = 15909
$ws.Range("M132").Value = -396423.14
$ws.Range("N132").Value = -20969
$ws.Range("H136").Value = 1240.5745
$ws.Range("I136").Value = 1333.7646
$ws.Range("J136").Value = 996.8461
$ws.Range("K136").Value = 4001.2938
$ws.Range("L136").Value = 2990.5383
$ws.Range("M136").Value = -1451.2938
$ws.Range("N136").Value = -8090.5383
$ws.Range("H139").Value = 42762.8
$ws.Range("J139").Value = 44180.89
$ws.Range("L139").Value = 44180.89
$ws.Range("N139").Value = -54460.89

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 48104.938
$ws.Range("I134").Value = 58794.348
$ws.Range("K134").Value = 176383.044
$ws.Range("M134").Value = -173848.044

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1361.2329
$ws.Range("I31").Value = 1074.3396
$ws.Range("J31").Value = 2121.5
$ws.Range("K31").Value = 1074.3396
$ws.Range("L31").Value = 2121.5
$ws.Range("M31").Value = -779.3396
$ws.Range("N31").Value = -2711.5
$ws.Range("H34").Value = 1361.2329
$ws.Range("I34").Value = 1074.3396
$ws.Range("J34").Value = 2121.5
$ws.Range("K34").Value = 1074.3396
$ws.Range("L34").Value = 2121.5
$ws.Range("M34").Value = -872.3396
$ws.Range("N34").Value = -2525.5
$ws.Range("H58").Value = 783.3684
$ws.Range("I58").Value = 851.6585
$ws.Range("J58").Value = 608.375
$ws.Range("K58").Value = 851.6585
$ws.Range("L58").Value = 608.375
$ws.Range("M58").Value = -648.6585
$ws.Range("N58").Value = -1014.375
$ws.Range("H132").Value = 2302.025
$ws.Range("I132").Value = 2083.818
$ws.Range("K132").Value = 6251.454000000001
$ws.Range("M132").Value = -3721.454000000001
$ws.Range("H136").Value = 783.3684
$ws.Range("I136").Value = 851.6585
$ws.Range("J136").Value = 608.375
$ws.Range("K136").Value = 2554.9755
$ws.Range("L136").Value = 1825.125
$ws.Range("M136").Value = -4.975500000000011
$ws.Range("N136").Value = -6925.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 35
$ws.Range("I23").Value = 25
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = 75
$ws.Range("L23").Value = 135
$ws.Range("M23").Value = 160
$ws.Range("N23").Value = -605
$ws.Range("H131").Value = 927.0405
$ws.Range("J131").Value = 982.18463
$ws.Range("L131").Value = 2946.55389
$ws.Range("N131").Value = -13026.55389

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1632.8549
$ws.Range("I132").Value = 1411.4
$ws.Range("J132").Value = 2035.5
$ws.Range("K132").Value = 4234.200000000001
$ws.Range("L132").Value = 6106.5
$ws.Range("M132").Value = -1704.200000000001
$ws.Range("N132").Value = -11166.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1720.44
$ws.Range("I132").Value = 1710.2285
$ws.Range("J132").Value = 1744.2667
$ws.Range("K132").Value = 5130.6855
$ws.Range("L132").Value = 5232.800099999999
$ws.Range("M132").Value = -2600.6855
$ws.Range("N132").Value = -10292.8001
$ws.Range("H136").Value = 1726.9131
$ws.Range("I136").Value = 1512.3125
$ws.Range("J136").Value = 2217.4285
$ws.Range("K136").Value = 4536.9375
$ws.Range("L136").Value = 6652.2855
$ws.Range("M136").Value = -1986.9375
$ws.Range("N136").Value = -11752.2855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 18391.25
$ws.Range("I75").Value = 5000
$ws.Range("J75").Value = 20304.285
$ws.Range("K75").Value = 5000
$ws.Range("L75").Value = 20304.285
$ws.Range("M75").Value = -4064
$ws.Range("N75").Value = -22176.285
$ws.Range("H78").Value = 18391.25
$ws.Range("I78").Value = 5000
$ws.Range("J78").Value = 20304.285
$ws.Range("K78").Value = 15000
$ws.Range("L78").Value = 60912.855
$ws.Range("N78").Value = -70272.855
